$d = $word.ActiveDocument

$d.Content.Find.Execute("2024-05-04 Saturday", $true, $false, $false, $false, $false, $true, 1, $false, "2024-05-05 Sunday", 2)
$d.Content.Find.Execute("810÷8=", $true, $false, $false, $false, $false, $true, 1, $false, "445÷3=", 2)
$d.Content.Find.Execute("482÷9=", $true, $false, $false, $false, $false, $true, 1, $false, "874÷6=", 2)
$d.Content.Find.Execute("382÷9=", $true, $false, $false, $false, $false, $true, 1, $false, "153÷6=", 2)
$d.Content.Find.Execute("894÷9=", $true, $false, $false, $false, $false, $true, 1, $false, "236÷8=", 2)
$d.Content.Find.Execute("145÷6=", $true, $false, $false, $false, $false, $true, 1, $false, "984÷2=", 2)
$d.Content.Find.Execute("250÷8=", $true, $false, $false, $false, $false, $true, 1, $false, "291÷4=", 2)
$d.Content.Find.Execute("826÷4=", $true, $false, $false, $false, $false, $true, 1, $false, "523÷8=", 2)
$d.Content.Find.Execute("922÷5=", $true, $false, $false, $false, $false, $true, 1, $false, "552÷5=", 2)
$d.Content.Find.Execute("497÷5=", $true, $false, $false, $false, $false, $true, 1, $false, "236÷9=", 2)
$d.Content.Find.Execute("913÷8=", $true, $false, $false, $false, $false, $true, 1, $false, "907÷6=", 2)
$d.Content.Find.Execute("108÷6=", $true, $false, $false, $false, $false, $true, 1, $false, "657÷2=", 2)
$d.Content.Find.Execute("673÷5=", $true, $false, $false, $false, $false, $true, 1, $false, "941÷8=", 2)
$d.Content.Find.Execute("361÷4=", $true, $false, $false, $false, $false, $true, 1, $false, "501÷9=", 2)
$d.Content.Find.Execute("600÷5=", $true, $false, $false, $false, $false, $true, 1, $false, "531÷9=", 2)
$d.Content.Find.Execute("418÷8=", $true, $false, $false, $false, $false, $true, 1, $false, "800÷7=", 2)
$d.Content.Find.Execute("757÷4=", $true, $false, $false, $false, $false, $true, 1, $false, "346÷8=", 2)
$d.Content.Find.Execute("194÷7=", $true, $false, $false, $false, $false, $true, 1, $false, "372÷7=", 2)
$d.Content.Find.Execute("857÷9=", $true, $false, $false, $false, $false, $true, 1, $false, "587÷8=", 2)
$d.Content.Find.Execute("199÷8=", $true, $false, $false, $false, $false, $true, 1, $false, "236÷3=", 2)
$d.Content.Find.Execute("930÷3=", $true, $false, $false, $false, $false, $true, 1, $false, "973÷8=", 2)
$d.Content.Find.Execute("993÷2=", $true, $false, $false, $false, $false, $true, 1, $false, "270÷8=", 2)
$d.Content.Find.Execute("478÷5=", $true, $false, $false, $false, $false, $true, 1, $false, "566÷3=", 2)
$d.Content.Find.Execute("177÷8=", $true, $false, $false, $false, $false, $true, 1, $false, "169÷3=", 2)
$d.Content.Find.Execute("374÷7=", $true, $false, $false, $false, $false, $true, 1, $false, "797÷8=", 2)
$d.Content.Find.Execute("834÷4=", $true, $false, $false, $false, $false, $true, 1, $false, "287÷5=", 2)
